# Update the "Förändrad" (Changed) date column (C) for every data row.
# The automatic update process advances the stored date serial by one day
# (45180 -> 45181, i.e. 2023-09-11 -> 2023-09-12) for all rows 2..118.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -ne $null) {
        $cell.Value2 = $cell.Value2 + 1
    }
}
